$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: clone formatting (borders/alignment/font) from matching template rows ---
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A23:E23").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A24:E24").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A25:E25").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A26:E26").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A27:E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 2: drop spurious cells the template paste introduced ---
$ws.Range("B18:E18").Clear() | Out-Null
$ws.Range("A25").Clear() | Out-Null
$ws.Range("A26").Clear() | Out-Null
$ws.Range("A27").Clear() | Out-Null

# --- Step 3: explicit row heights ---
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 43.2
$ws.Rows.Item(20).RowHeight = 43.2
$ws.Rows.Item(21).RowHeight = 43.2
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 21.6
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 21.6
$ws.Rows.Item(26).RowHeight = 21.6
$ws.Rows.Item(27).RowHeight = 31.8

# --- Step 4: cell values ---
$ws.Range("A16").Value = 'SCRIPT/T01P01A/um1201.ssb'
$ws.Range("B16").Value = 182
$ws.Range("C16").Value = ' Hi! Isn\''t the weather wonderful?'
$ws.Range("D16").Value = ' Привет! Сегодня чудесная\nпогода, да?'
$ws.Range("E16").Value = ' Ðñéâåó! Òåãïäîÿ œôäåòîàÿ\nðïãïäà, äà?'
$ws.Range("A17").Value = 'SCRIPT/T01P01A/um1301.ssb'
$ws.Range("B17").Value = 163
$ws.Range("C17").Value = ' The guild\''s Pokémon seem to be\nin a panic or something…'
$ws.Range("D17").Value = ' Все Покемоны гильдии в какой-то\nпанике или что-то вроде того...'
$ws.Range("E17").Value = ' Âòå Ðïëåíïîú ãéìûäéé â ëàëïê-óï\nðàîéëå éìé œóï-óï âñïäå óïãï…'
$ws.Range("A18").Value = 'SCRIPT/T01P01A/um1310.ssb'
$ws.Range("A19").Value = 'SCRIPT/T01P01A/um1313.ssb'
$ws.Range("B19").ClearContents() | Out-Null
$ws.Range("C19").ClearContents() | Out-Null
$ws.Range("D19").ClearContents() | Out-Null
$ws.Range("E19").ClearContents() | Out-Null
$ws.Range("A20").Value = 'SCRIPT/T01P01A/um1316.ssb'
$ws.Range("B20").Value = 141
$ws.Range("C20").Value = ' I\''ve heard.[K] Your whole guild is\ntrying to capture [CS:N]Grovyle[CR]?'
$ws.Range("D20").Value = ' Я всё слышал.[K] Вся ваша гильдия\nпытается поймать [CS:N]Гровайла[CR]?'
$ws.Range("E20").Value = ' Ÿ âòæ òìúšàì.[K] Âòÿ âàšà ãéìûäéÿ\nðúóàåóòÿ ðïêíàóû [CS:N]Ãñïâàêìà[CR]?'
$ws.Range("A21").Value = 'SCRIPT/T01P01A/um1401.ssb'
$ws.Range("B21").Value = 144
$ws.Range("C21").Value = ' I hope you succeed!'
$ws.Range("D21").Value = ' Надеюсь, у вас всё получится!'
$ws.Range("E21").Value = ' Îàäåýòû, ô âàò âòæ ðïìôœéóòÿ!'
$ws.Range("A22").Value = 'SCRIPT/T01P01A/um1601.ssb'
$ws.Range("B22").Value = 119
$ws.Range("C22").Value = ' He\''s trying to paralyze the\nplanet by stealing Time Gears, isn\''t he?'
$ws.Range("D22").Value = ' Он пытается парализовать\nпланету похищая Шестерни Времени, так?'
$ws.Range("E22").Value = ' Ïî ðúóàåóòÿ ðàñàìéèïâàóû\nðìàîåóô ðïöéþàÿ Šåòóåñîé Âñåíåîé, óàë?'
$ws.Range("B23").Value = 122
$ws.Range("C23").Value = ' What a horrible thing to do!\nThat [CS:N]Grovyle[CR]\''s horrid!'
$ws.Range("D23").Value = ' Какой ужасный поступок! Какой\nужасный [CS:N]Гровайл[CR]!'
$ws.Range("E23").Value = ' Ëàëïê ôçàòîúê ðïòóôðïë! Ëàëïê\nôçàòîúê [CS:N]Ãñïâàêì[CR]!'
$ws.Range("A24").Value = 'SCRIPT/T01P01A/um1607.ssb'
$ws.Range("B24").Value = 91
$ws.Range("C24").Value = ' I think we\''re in safe hands\nwith the great [CS:N]Dusknoir[CR].'
$ws.Range("D24").Value = ' Я думаю, что с великим\n[CS:N]Даскнуаром[CR] мы в безопасности.'
$ws.Range("E24").Value = ' Ÿ äôíàý, œóï ò âåìéëéí\n[CS:N]Äàòëîôàñïí[CR] íú â áåèïðàòîïòóé.'
$ws.Range("B25").Value = 94
$ws.Range("C25").Value = ' I\''m sure he\''ll catch [CS:N]Grovyle[CR]\nfor us!'
$ws.Range("D25").Value = ' Я уверен, что он поймает\n[CS:N]Гровайла[CR]!'
$ws.Range("E25").Value = ' Ÿ ôâåñåî, œóï ïî ðïêíàåó\n[CS:N]Ãñïâàêìà[CR]!'
$ws.Range("B26").Value = 97
$ws.Range("C26").Value = ' ...Oops. I shouldn\''t have said\nthat so loud.'
$ws.Range("D26").Value = ' ...Упс. Я не должен был это\nвыкрикивать.'
$ws.Range("E26").Value = ' ...Ôðò. Ÿ îå äïìçåî áúì üóï\nâúëñéëéâàóû.'
$ws.Range("B27").Value = 100
$ws.Range("C27").Value = ' Who knows? [CS:N]Grovyle[CR] might\neven be listening in on us.[K] One can\''t be\ntoo careful.'
$ws.Range("D27").Value = ' Кто знает? Может быть [CS:N]Гровайл[CR]\nнас подслушивает.[K] Нужно быть осторожнее.'
$ws.Range("E27").Value = ' Ëóï èîàåó? Íïçåó áúóû [CS:N]Ãñïâàêì[CR]\nîàò ðïäòìôšéâàåó.[K] Îôçîï áúóû ïòóïñïçîåå.'

# --- Step 5: dimension / view state ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("E27").Select() | Out-Null

Write-Host "done"
